$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 6) - data now only spans rows 1-5
$ws.Rows("6").Delete()

# New data values (TPM update) for rows 2-5, columns A-T
$data = @(
    @("ECs",           "Qrfp", "Qrfpr", "FAPs", 1, 0.5,                0.259189,            0.518378, 0.1837975229301632, 0.1558987974170948, 1, 0.3333333333333333, 0.06285, 0.18855, 1, 1, 0.01629002865,  0.09774017189999999, 0.1837975229301632, 0.1558987974170948),
    @("FAPs",          "Qrfp", "Qrfpr", "FAPs", 2, 0.6666666666666666, 0.3205933333333333, 0.9617800000000001, 0.2273409000173263, 0.2892490525828902, 1, 0.3333333333333333, 0.06285, 0.18855, 1, 1, 0.020149291,    0.181343619,          0.2273409000173263, 0.2892490525828902),
    @("MuSCs",         "Qrfp", "Qrfpr", "FAPs", 2, 1,                  0.646281,            1.292562, 0.458294321583203,  0.3887295783907397, 1, 0.3333333333333333, 0.06285, 0.18855, 1, 1, 0.04061876085,  0.2437125651,         0.458294321583203,  0.3887295783907397),
    @("Resolving-Mac", "Qrfp", "Qrfpr", "FAPs", 2, 0.6666666666666666, 0.1841243333333333, 0.552373, 0.1305672554693075, 0.1661225716092753, 1, 0.3333333333333333, 0.06285, 0.18855, 1, 1, 0.01157221435,  0.10414992915,        0.1305672554693075, 0.1661225716092753)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
